$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1315.5
$ws.Range("I28").Value = 1315.5
$ws.Range("K28").Value = 1315.5
$ws.Range("M28").Value = -830.5
$ws.Range("H38").Value = 432.8125
$ws.Range("I38").Value = 350.5
$ws.Range("K38").Value = 1051.5
$ws.Range("M38").Value = -679.5
$ws.Range("H92").Value = 395.8889
$ws.Range("I92").Value = 366.2353
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 366.2353
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 881.7646999999999
$ws.Range("N92").Value = -3396
$ws.Range("H99").Value = 200200220
$ws.Range("I99").Value = 380
$ws.Range("K99").Value = 1140
$ws.Range("M99").Value = 358
$ws.Range("H127").Value = 115192.125
$ws.Range("I127").Value = 131519.58
$ws.Range("K127").Value = 394558.74
$ws.Range("M127").Value = -389598.74
$ws.Range("H129").Value = 1025.1111
$ws.Range("I129").Value = 596.1667
$ws.Range("J129").Value = 1883
$ws.Range("K129").Value = 1788.5001
$ws.Range("L129").Value = 5649
$ws.Range("M129").Value = 3211.4999
$ws.Range("N129").Value = -15649
$ws.Range("H132").Value = 7209.8125
$ws.Range("I132").Value = 2115
$ws.Range("J132").Value = 13760.286
$ws.Range("K132").Value = 6345
$ws.Range("L132").Value = 41280.858
$ws.Range("M132").Value = -3815
$ws.Range("N132").Value = -46340.858
$ws.Range("H137").Value = 33692.516
$ws.Range("I137").Value = 1422.5769
$ws.Range("J137").Value = 201496.2
$ws.Range("K137").Value = 4267.7307
$ws.Range("L137").Value = 604488.6000000001
$ws.Range("M137").Value = -1717.7307
$ws.Range("N137").Value = -609588.6000000001
$ws.Range("H141").Value = 67285.57000000001
$ws.Range("I141").Value = 92979.39999999999
$ws.Range("K141").Value = 278938.2
$ws.Range("M141").Value = -273758.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2814.7144
$ws.Range("I2").Value = 2545.5
$ws.Range("K2").Value = 2545.5
$ws.Range("M2").Value = -2432.5
$ws.Range("H45").Value = 25868.615
$ws.Range("I45").Value = 30278.572
$ws.Range("K45").Value = 30278.572
$ws.Range("M45").Value = -29901.572
$ws.Range("H61").Value = 3055.3076
$ws.Range("I61").Value = 2268.5
$ws.Range("J61").Value = 3729.7144
$ws.Range("K61").Value = 2268.5
$ws.Range("L61").Value = 3729.7144
$ws.Range("M61").Value = -2056.5
$ws.Range("N61").Value = -4153.7144
$ws.Range("H74").Value = 1420.6666
$ws.Range("I74").Value = 1420.6666
$ws.Range("K74").Value = 1420.6666
$ws.Range("M74").Value = -546.6666
$ws.Range("H77").Value = 1420.6666
$ws.Range("I77").Value = 1420.6666
$ws.Range("K77").Value = 7103.333000000001
$ws.Range("M77").Value = -2735.333000000001
$ws.Range("H116").Value = 2814.7144
$ws.Range("I116").Value = 2545.5
$ws.Range("K116").Value = 2545.5
$ws.Range("M116").Value = -251.5
$ws.Range("H122").Value = 12535.632
$ws.Range("I122").Value = 12954.277
$ws.Range("K122").Value = 38862.831
$ws.Range("M122").Value = -36412.831
$ws.Range("H136").Value = 3055.3076
$ws.Range("I136").Value = 2268.5
$ws.Range("J136").Value = 3729.7144
$ws.Range("K136").Value = 6805.5
$ws.Range("L136").Value = 11189.1432
$ws.Range("M136").Value = -4255.5
$ws.Range("N136").Value = -16289.1432

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2814.7144
$ws.Range("I3").Value = 2545.5
$ws.Range("K3").Value = 2545.5
$ws.Range("M3").Value = -2431.5
$ws.Range("H107").Value = 19935.2
$ws.Range("I107").Value = 6521.7144
$ws.Range("J107").Value = 51233.332
$ws.Range("K107").Value = 6521.7144
$ws.Range("L107").Value = 51233.332
$ws.Range("M107").Value = -4601.7144
$ws.Range("N107").Value = -55073.332
$ws.Range("H134").Value = 2949.238
$ws.Range("I134").Value = 2572.6155
$ws.Range("K134").Value = 7717.8465
$ws.Range("M134").Value = -5182.8465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2098
$ws.Range("I16").Value = 1784.625
$ws.Range("K16").Value = 1784.625
$ws.Range("M16").Value = -1497.625
$ws.Range("H58").Value = 1678.6842
$ws.Range("I58").Value = 1460
$ws.Range("K58").Value = 1460
$ws.Range("M58").Value = -1257
$ws.Range("H113").Value = 2098
$ws.Range("I113").Value = 1784.625
$ws.Range("K113").Value = 1784.625
$ws.Range("M113").Value = 385.375
$ws.Range("H132").Value = 2404.3635
$ws.Range("I132").Value = 2423.8572
$ws.Range("K132").Value = 7271.571599999999
$ws.Range("M132").Value = -4741.571599999999
$ws.Range("H134").Value = 2066.6333
$ws.Range("I134").Value = 1927.16
$ws.Range("K134").Value = 5781.48
$ws.Range("M134").Value = -3246.48
$ws.Range("H136").Value = 1678.6842
$ws.Range("I136").Value = 1460
$ws.Range("K136").Value = 4380
$ws.Range("M136").Value = -1830

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1291.6
$ws.Range("J68").Value = 1616.5714
$ws.Range("L68").Value = 4849.7142
$ws.Range("N68").Value = -6471.7142
$ws.Range("H71").Value = 1291.6
$ws.Range("J71").Value = 1616.5714
$ws.Range("L71").Value = 14549.1426
$ws.Range("N71").Value = -22661.1426
$ws.Range("H81").Value = 55565444
$ws.Range("I81").Value = 1747
$ws.Range("J81").Value = 71440780
$ws.Range("K81").Value = 5241
$ws.Range("L81").Value = 214322340
$ws.Range("M81").Value = -4118
$ws.Range("N81").Value = -214324586
$ws.Range("H84").Value = 55565444
$ws.Range("I84").Value = 1747
$ws.Range("J84").Value = 71440780
$ws.Range("K84").Value = 15723
$ws.Range("L84").Value = 642967020
$ws.Range("M84").Value = -10107
$ws.Range("N84").Value = -642978252

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1505000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1505000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1505000
$ws.Range("N7").Value = -1505224
$ws.Range("H8").Value = 1505000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1505000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1505000
$ws.Range("N8").Value = -1505278
$ws.Range("H10").Value = 1005000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10338
$ws.Range("H102").Value = 1439.091
$ws.Range("I102").Value = 832.8570999999999
$ws.Range("K102").Value = 832.8570999999999
$ws.Range("M102").Value = 789.1429000000001
$ws.Range("H113").Value = 71431470
$ws.Range("I113").Value = 83335710
$ws.Range("K113").Value = 83335710
$ws.Range("M113").Value = -83333540
$ws.Range("H122").Value = 1620
$ws.Range("I122").Value = 1525
$ws.Range("K122").Value = 4575
$ws.Range("M122").Value = -2125
$ws.Range("H126").Value = 3173.5
$ws.Range("I126").Value = 2333
$ws.Range("J126").Value = 4014
$ws.Range("K126").Value = 6999
$ws.Range("L126").Value = 12042
$ws.Range("M126").Value = -4529
$ws.Range("N126").Value = -16982
$ws.Range("M7").ClearContents()
$ws.Range("M8").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6394.1665
$ws.Range("I40").Value = 6385
$ws.Range("K40").Value = 6385
$ws.Range("M40").Value = -6249
$ws.Range("H46").Value = 2711.4666
$ws.Range("I46").Value = 2257
$ws.Range("J46").Value = 3109.125
$ws.Range("K46").Value = 2257
$ws.Range("L46").Value = 3109.125
$ws.Range("M46").Value = -2069
$ws.Range("N46").Value = -3485.125
$ws.Range("H61").Value = 12889
$ws.Range("I61").Value = 16248.75
$ws.Range("K61").Value = 16248.75
$ws.Range("M61").Value = -16046.75
$ws.Range("H93").Value = 2198
$ws.Range("I93").Value = 1997.5
$ws.Range("K93").Value = 1997.5
$ws.Range("M93").Value = -749.5
$ws.Range("H100").Value = 4374.75
$ws.Range("I100").Value = 2999.3333
$ws.Range("K100").Value = 2999.3333
$ws.Range("M100").Value = -2458.3333
$ws.Range("H106").Value = 29563.334
$ws.Range("J106").Value = 29563.334
$ws.Range("L106").Value = 29563.334
$ws.Range("N106").Value = -32087.334
$ws.Range("H113").Value = 12889
$ws.Range("I113").Value = 16248.75
$ws.Range("K113").Value = 16248.75
$ws.Range("M113").Value = -14078.75
$ws.Range("H122").Value = 4633
$ws.Range("I122").Value = 3949.5
$ws.Range("K122").Value = 11848.5
$ws.Range("M122").Value = -9398.5
$ws.Range("H132").Value = 2926.0967
$ws.Range("I132").Value = 2760.2273
$ws.Range("J132").Value = 3331.5557
$ws.Range("K132").Value = 8280.6819
$ws.Range("L132").Value = 9994.667099999999
$ws.Range("M132").Value = -5750.6819
$ws.Range("N132").Value = -15054.6671
$ws.Range("H136").Value = 3249.1333
$ws.Range("I136").Value = 2374
$ws.Range("J136").Value = 3832.5557
$ws.Range("K136").Value = 7122
$ws.Range("L136").Value = 11497.6671
$ws.Range("M136").Value = -4572
$ws.Range("N136").Value = -16597.6671

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 20077.5
$ws.Range("J70").Value = 20077.5
$ws.Range("L70").Value = 20077.5
$ws.Range("N70").Value = -20707.5
$ws.Range("H73").Value = 20077.5
$ws.Range("J73").Value = 20077.5
$ws.Range("L73").Value = 20077.5
$ws.Range("N73").Value = -22261.5
$ws.Range("H100").Value = 4539.909
$ws.Range("I100").Value = 5740
$ws.Range("K100").Value = 11480
$ws.Range("M100").Value = -10939
$ws.Range("H124").Value = 53720.57
$ws.Range("J124").Value = 53720.57
$ws.Range("L124").Value = 53720.57
$ws.Range("N124").Value = -63540.57
$ws.Range("H126").Value = 2539.4443
$ws.Range("I126").Value = 1969
$ws.Range("J126").Value = 2995.8
$ws.Range("K126").Value = 5907
$ws.Range("L126").Value = 8987.400000000001
$ws.Range("M126").Value = -3437
$ws.Range("N126").Value = -13927.4
$ws.Range("H132").Value = 7400.7354
$ws.Range("I132").Value = 7737.3667
$ws.Range("K132").Value = 23212.1001
$ws.Range("M132").Value = -20682.1001
$ws.Range("H136").Value = 3710.5833
$ws.Range("I136").Value = 2452.8
$ws.Range("J136").Value = 9999.5
$ws.Range("K136").Value = 7358.400000000001
$ws.Range("L136").Value = 29998.5
$ws.Range("M136").Value = -4808.400000000001
$ws.Range("N136").Value = -35098.5
